$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 432; this shifts the existing
# rows 432..451 down to 433..452 (old row 451 -> new row 452) and leaves
# an empty row 432 for the new weekly record.
$ws.Rows.Item(432).Insert()

# Populate the newly inserted row 432 with the new weekly price record.
$ws.Range("A432").Value = 4
$ws.Range("B432").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C432").Value = "Los Lagos"
$ws.Range("D432").Value = 45041
$ws.Range("D432").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E432").Value = 10
$ws.Range("F432").Value = 100112017
$ws.Range("G432").Value = "Apio"
$ws.Range("H432").Value = "Americana (o)"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 45
$ws.Range("K432").Value = 12000
$ws.Range("L432").Value = 12000
$ws.Range("M432").Value = 12000
$ws.Range("N432").Value = "$/docena de matas"
$ws.Range("O432").Value = "Región de Coquimbo"
$ws.Range("P432").Value = 2000
$ws.Range("Q432").Value = 6
$ws.Range("R432").Value = "Hortaliza"
